# Updates cryptos list values (price + 1h volume change) per upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''27.070.18'
$ws.Range("D3").Value = '''1.566.93'
$ws.Range("E3").Value = '  +1.06%  '
$ws.Range("E4").Value = '  +0.53%  '
$ws.Range("D5").Value = '''208.61'
$ws.Range("E5").Value = '  +1.10%  '
$ws.Range("E6").Value = '  +0.75%  '
$ws.Range("E7").Value = '  +0.52%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("E9").Value = '  +1.07%  '
$ws.Range("E10").Value = '  +1.71%  '
$ws.Range("D11").Value = '''0.0861'
$ws.Range("E11").Value = '  +0.66%  '
$ws.Range("D12").Value = '''1.571.40'
$ws.Range("E12").Value = '  +1.40%  '
$ws.Range("D13").Value = '''3.78'
$ws.Range("E13").Value = '  +1.23%  '
$ws.Range("E14").Value = '  +0.29%  '
$ws.Range("D15").Value = '''27.062.35'
$ws.Range("E15").Value = '  +0.60%  '
$ws.Range("D16").Value = '''61.91'
$ws.Range("E16").Value = '  +0.40%  '
$ws.Range("E17").Value = '  +0.97%  '
$ws.Range("D18").Value = '''7.42'
$ws.Range("E18").Value = '  +2.11%  '
$ws.Range("D19").Value = '''215.71'
$ws.Range("E19").Value = '  -0.66%  '
$ws.Range("E20").Value = '  +0.53%  '
$ws.Range("E21").Value = '  +2.29%  '
$ws.Range("D22").Value = '''9.21'
$ws.Range("E22").Value = '  -0.16%  '
$ws.Range("D24").Value = '''154.11'
$ws.Range("E24").Value = '  +0.14%  '
$ws.Range("E25").Value = '  -0.30%  '
$ws.Range("D26").Value = '''15.04'
$ws.Range("E26").Value = '  +0.50%  '
$ws.Range("D27").Value = '''0.106'
$ws.Range("E27").Value = '  +1.46%  '
$ws.Range("E28").Value = '  +0.61%  '
$ws.Range("E29").Value = '  +1.47%  '
$ws.Range("E30").Value = '  +3.93%  '
$ws.Range("E31").Value = '  +0.55%  '
$ws.Range("E32").Value = '  +3.69%  '
$ws.Range("D33").Value = '''1.423.68'
$ws.Range("E33").Value = '  +0.69%  '
$ws.Range("D34").Value = '''1.09'
$ws.Range("E34").Value = '  +13.03%  '
$ws.Range("E35").Value = '  +1.41%  '
$ws.Range("D36").Value = '''2.35'
$ws.Range("E36").Value = '  +2.87%  '
$ws.Range("D37").Value = '''0.0167'
$ws.Range("E37").Value = '  +0.95%  '
$ws.Range("E38").Value = '  +1.30%  '
$ws.Range("E39").Value = '  +2.32%  '
$ws.Range("D40").Value = '''0.812'
$ws.Range("E40").Value = '  +0.49%  '
$ws.Range("E41").Value = '  +0.55%  '
$ws.Range("D42").Value = '''2.35'
$ws.Range("E42").Value = '  +1.19%  '
$ws.Range("E43").Value = '  +0.39%  '
$ws.Range("D44").Value = '''64.73'
$ws.Range("E44").Value = '  +0.08%  '
$ws.Range("E45").Value = '  -0.40%  '
$ws.Range("D46").Value = '''1.703.00'
$ws.Range("E46").Value = '  +1.04%  '
$ws.Range("D47").Value = '''86.76'
$ws.Range("E47").Value = '  -0.85%  '
$ws.Range("D48").Value = '''0.0₆0103'
$ws.Range("E48").Value = '  +2.91%  '
$ws.Range("D49").Value = '''0.0518'
$ws.Range("E49").Value = '  +0.47%  '
$ws.Range("D50").Value = '''0.0963'
$ws.Range("E50").Value = '  +0.21%  '
$ws.Range("E51").Value = '  +0.47%  '
